# Insert a new first sheet "2023-10-04" ahead of the existing "2023-09-01" sheet,
# containing the header row plus a single data row (the latest resume extraction).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "2023-10-04"

$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Mobile No."
$ws.Range("C1").Value = "Skills"

# Mobile numbers are stored as text in this workbook (leading "+" / leading
# zeros elsewhere would otherwise be mangled), so force the column to Text
# before writing the digits-only number.
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = "maheshkulkarni01121@gmail.com"
$ws.Range("B2").Value = "9423627124"
$ws.Range("C2").Value = "Editing, Programming, English, Python, Css, Content, Training, Networking, Research, Database, Engineering, Php, Website, Analysis, Tensorflow, C, Html, Technical, Writing, Video"

$ws.Activate()
